$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.306.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "'1.868.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'235.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4683"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.2843"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").Value = "'0.06532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").Value = "'21.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.88%  "
$ws.Range("D11").Value = "'0.07867"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").Value = "'97.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "'1.864.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "'5.100"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'0.6757"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'279.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "'30.304.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'1.0000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'5.504"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").Value = "'12.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "'2.111.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'0.000007297"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.162"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'165.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.174"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "'19.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "'1.931"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("D29").Value = "'1.376"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "'0.09634"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'4.365"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "'1.478"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "'4.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'0.04705"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "'1.128"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("D36").Value = "'0.7071"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "'2.721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'0.01854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'6.279"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").Value = "'2.531"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'73.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").Value = "'1.946"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "'0.8496"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "'0.4178"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'7.166"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "'9.278"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "'935.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.05%  "
$ws.Range("D50").Value = "'34.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").Value = "'0.1124"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
